# "Updated xlsx with correct color coding and tokens"
#
# 1) The gpt-35-turbo region (H5) had its throttled token quota corrected
#    from 240000 down to 40000 (dependent formulas I5/J5 recalc automatically).
# 2) The newly-added gpt-35-turbo-16k row (row 6) was left with the default/
#    unstyled formatting when it was first added; this re-applies the same
#    color coding (fill/border/number format) used by the other model rows
#    (copied here from row 5, which shares row 6's color group) so the whole
#    table is visually consistent again.
# 3) Leave the selection parked on H6, matching where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the throttled token quota for gpt-35-turbo (South Central) ---
$ws.Range("H5").Value = 40000

# --- Re-apply correct color coding to the gpt-35-turbo-16k row ---
$ws.Range("A5:J5").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Park the selection on H6 ---
$ws.Range("H6").Select() | Out-Null
